$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 519.369
$ws.Range("I15").Value = 519.369
$ws.Range("K15").Value = 1558.107
$ws.Range("M15").Value = -1389.107

# ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2279.5
$ws.Range("J121").Value = 2279.5
$ws.Range("L121").Value = 6838.5
$ws.Range("N121").Value = -10332.5

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 7610.4585
$ws.Range("I135").Value = 7523
$ws.Range("J135").Value = 7662.933
$ws.Range("K135").Value = 67707
$ws.Range("L135").Value = 68966.397
$ws.Range("M135").Value = -65172
$ws.Range("N135").Value = -74036.397

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7629.3677
$ws.Range("I137").Value = 3224.56
$ws.Range("K137").Value = 9673.68
$ws.Range("M137").Value = -7123.68

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5245.4434
$ws.Range("I138").Value = 5324.6113
$ws.Range("J138").Value = 5227.4053
$ws.Range("K138").Value = 15973.8339
$ws.Range("L138").Value = 15682.2159
$ws.Range("M138").Value = -10833.8339
$ws.Range("N138").Value = -25962.2159

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4232.1763
$ws.Range("I2").Value = 1099.75
$ws.Range("J2").Value = 5196
$ws.Range("K2").Value = 1099.75
$ws.Range("L2").Value = 5196
$ws.Range("M2").Value = -986.75
$ws.Range("N2").Value = -5422

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1560997.5
$ws.Range("I32").Value = 4556884
$ws.Range("K32").Value = 4556884
$ws.Range("M32").Value = -4556597

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11813.885
$ws.Range("I61").Value = 2461.7
$ws.Range("J61").Value = 42987.832
$ws.Range("K61").Value = 2461.7
$ws.Range("L61").Value = 42987.832
$ws.Range("M61").Value = -2249.7
$ws.Range("N61").Value = -43411.832

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 25468.5
$ws.Range("I74").Value = 3469.6667
$ws.Range("J74").Value = 29868.268
$ws.Range("K74").Value = 3469.6667
$ws.Range("L74").Value = 29868.268
$ws.Range("M74").Value = -2595.6667
$ws.Range("N74").Value = -31616.268

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 25468.5
$ws.Range("I77").Value = 3469.6667
$ws.Range("J77").Value = 29868.268
$ws.Range("K77").Value = 17348.3335
$ws.Range("L77").Value = 149341.34
$ws.Range("M77").Value = -12980.3335
$ws.Range("N77").Value = -158077.34

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 14535.117
$ws.Range("I102").Value = 2473.2
$ws.Range("K102").Value = 2473.2
$ws.Range("M102").Value = -851.1999999999998

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 4232.1763
$ws.Range("I116").Value = 1099.75
$ws.Range("J116").Value = 5196
$ws.Range("K116").Value = 1099.75
$ws.Range("L116").Value = 5196
$ws.Range("M116").Value = 1194.25
$ws.Range("N116").Value = -9784

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4779208
$ws.Range("I132").Value = 8371.666999999999
$ws.Range("K132").Value = 25115.001
$ws.Range("M132").Value = -22585.001

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 11813.885
$ws.Range("I136").Value = 2461.7
$ws.Range("J136").Value = 42987.832
$ws.Range("K136").Value = 7385.099999999999
$ws.Range("L136").Value = 128963.496
$ws.Range("M136").Value = -4835.099999999999
$ws.Range("N136").Value = -134063.496

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4232.1763
$ws.Range("I3").Value = 1099.75
$ws.Range("J3").Value = 5196
$ws.Range("K3").Value = 1099.75
$ws.Range("L3").Value = 5196
$ws.Range("M3").Value = -985.75
$ws.Range("N3").Value = -5424

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3999.6553
$ws.Range("I86").Value = 3878.7222
$ws.Range("J86").Value = 4197.5454
$ws.Range("K86").Value = 3878.7222
$ws.Range("L86").Value = 4197.5454
$ws.Range("M86").Value = -2755.7222
$ws.Range("N86").Value = -6443.5454

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3999.6553
$ws.Range("I89").Value = 3878.7222
$ws.Range("J89").Value = 4197.5454
$ws.Range("K89").Value = 19393.611
$ws.Range("L89").Value = 20987.727
$ws.Range("M89").Value = -13777.611
$ws.Range("N89").Value = -32219.727

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4805.35
$ws.Range("I107").Value = 6458.615
$ws.Range("J107").Value = 1735
$ws.Range("K107").Value = 6458.615
$ws.Range("L107").Value = 1735
$ws.Range("M107").Value = -4538.615
$ws.Range("N107").Value = -5575

# BSM row 108
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 175000
$ws.Range("I108").Value = 200000
$ws.Range("J108").Value = 150000
$ws.Range("K108").Value = 200000
$ws.Range("L108").Value = 150000
$ws.Range("M108").Value = -196160
$ws.Range("N108").Value = -157680

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8148.909
$ws.Range("I134").Value = 2513.4614
$ws.Range("K134").Value = 7540.3842
$ws.Range("M134").Value = -5005.3842

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3503.9
$ws.Range("I22").Value = 3907.8
$ws.Range("J22").Value = 3100
$ws.Range("K22").Value = 3907.8
$ws.Range("L22").Value = 3100
$ws.Range("M22").Value = -3557.8
$ws.Range("N22").Value = -3800

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21119.139
$ws.Range("I31").Value = 7180.684
$ws.Range("J31").Value = 47602.2
$ws.Range("K31").Value = 7180.684
$ws.Range("L31").Value = 47602.2
$ws.Range("M31").Value = -6885.684
$ws.Range("N31").Value = -48192.2

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 21119.139
$ws.Range("I34").Value = 7180.684
$ws.Range("J34").Value = 47602.2
$ws.Range("K34").Value = 7180.684
$ws.Range("L34").Value = 47602.2
$ws.Range("M34").Value = -6978.684
$ws.Range("N34").Value = -48006.2

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5926.1143
$ws.Range("I132").Value = 1667.7727
$ws.Range("K132").Value = 5003.3181
$ws.Range("M132").Value = -2473.3181

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 19612124
$ws.Range("I134").Value = 1266
$ws.Range("K134").Value = 3798
$ws.Range("M134").Value = -1263

# CUL row 36
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 3970.875
$ws.Range("I36").Value = 2558.4
$ws.Range("J36").Value = 6325
$ws.Range("K36").Value = 7675.200000000001
$ws.Range("L36").Value = 18975
$ws.Range("M36").Value = -7506.200000000001
$ws.Range("N36").Value = -19313

# CUL row 114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2584.2727
$ws.Range("J114").Value = 2807.625
$ws.Range("L114").Value = 8422.875
$ws.Range("N114").Value = -14930.875

# GSM row 26
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 41797.7
$ws.Range("I26").Value = 19998.5
$ws.Range("J26").Value = 47247.5
$ws.Range("K26").Value = 19998.5
$ws.Range("L26").Value = 47247.5
$ws.Range("M26").Value = -19718.5
$ws.Range("N26").Value = -47807.5

# GSM row 50
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 41797.7
$ws.Range("I50").Value = 19998.5
$ws.Range("J50").Value = 47247.5
$ws.Range("K50").Value = 19998.5
$ws.Range("L50").Value = 47247.5
$ws.Range("M50").Value = -19500.5
$ws.Range("N50").Value = -48243.5

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 102160.445
$ws.Range("I113").Value = 181090.6
$ws.Range("K113").Value = 181090.6
$ws.Range("M113").Value = -178920.6

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 76489.47
$ws.Range("I132").Value = 94666.17999999999
$ws.Range("J132").Value = 26503.5
$ws.Range("K132").Value = 283998.54
$ws.Range("L132").Value = 79510.5
$ws.Range("M132").Value = -281468.54
$ws.Range("N132").Value = -84570.5

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2749.2942
$ws.Range("I16").Value = 3076.9092
$ws.Range("K16").Value = 3076.9092
$ws.Range("M16").Value = -2906.9092

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11353.071
$ws.Range("I22").Value = 2358.3333
$ws.Range("K22").Value = 2358.3333
$ws.Range("M22").Value = -2063.3333

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 11353.071
$ws.Range("I27").Value = 2358.3333
$ws.Range("K27").Value = 2358.3333
$ws.Range("M27").Value = -2251.3333

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1669216.5
$ws.Range("I46").Value = 3333933.8
$ws.Range("K46").Value = 3333933.8
$ws.Range("M46").Value = -3333745.8

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7208.3
$ws.Range("I68").Value = 2260.5
$ws.Range("J68").Value = 26999.5
$ws.Range("K68").Value = 2260.5
$ws.Range("L68").Value = 26999.5
$ws.Range("M68").Value = -1511.5
$ws.Range("N68").Value = -28497.5

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 7208.3
$ws.Range("I71").Value = 2260.5
$ws.Range("J71").Value = 26999.5
$ws.Range("K71").Value = 11302.5
$ws.Range("L71").Value = 134997.5
$ws.Range("M71").Value = -7558.5
$ws.Range("N71").Value = -142485.5

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3097.6
$ws.Range("J81").Value = 3249.5
$ws.Range("L81").Value = 6499
$ws.Range("N81").Value = -8621

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3097.6
$ws.Range("J84").Value = 3249.5
$ws.Range("L84").Value = 32495
$ws.Range("N84").Value = -43103

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3952.95
$ws.Range("J113").Value = 7809.8887
$ws.Range("L113").Value = 23429.6661
$ws.Range("N113").Value = -27769.6661

